# Auto-generated edit script: updates Price (D), Volume(1h) (E) and Hora (G) columns
# for the "cryptos" symbol-list refresh committed on Thu Dec 29 20:15:34 UTC 2022.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.06"
$ws.Range("G2").Value = "'20"
$ws.Range("D3").Value = "'24.22"
$ws.Range("G3").Value = "'20"
$ws.Range("D4").Value = "'5.288"
$ws.Range("G4").Value = "'20"
$ws.Range("G5").Value = "'20"
$ws.Range("D6").Value = "'6.490"
$ws.Range("G6").Value = "'20"
$ws.Range("D7").Value = "'3.146"
$ws.Range("G7").Value = "'20"
$ws.Range("D8").Value = "'0.8174"
$ws.Range("G8").Value = "'20"
$ws.Range("D9").Value = "'0.8582"
$ws.Range("G9").Value = "'20"
$ws.Range("D10").Value = "'0.1366"
$ws.Range("G10").Value = "'20"
$ws.Range("D11").Value = "'0.06959"
$ws.Range("G11").Value = "'20"
$ws.Range("D12").Value = "'0.03142"
$ws.Range("G12").Value = "'20"
$ws.Range("G13").Value = "'20"
$ws.Range("D14").Value = "'0.09403"
$ws.Range("G14").Value = "'20"
$ws.Range("D15").Value = "'3.783"
$ws.Range("G15").Value = "'20"
$ws.Range("G16").Value = "'20"
$ws.Range("D17").Value = "'0.04682"
$ws.Range("G17").Value = "'20"
$ws.Range("D18").Value = "'0.0005994"
$ws.Range("G18").Value = "'20"
$ws.Range("D19").Value = "'0.006168"
$ws.Range("G19").Value = "'20"
$ws.Range("D20").Value = "'0.001240"
$ws.Range("G20").Value = "'20"
$ws.Range("D21").Value = "'0.004618"
$ws.Range("G21").Value = "'20"
$ws.Range("D22").Value = "'0.00006104"
$ws.Range("G22").Value = "'20"
$ws.Range("D23").Value = "'3.499"
$ws.Range("G23").Value = "'20"
$ws.Range("D24").Value = "'2.147"
$ws.Range("G24").Value = "'20"
$ws.Range("D25").Value = "'0.3197"
$ws.Range("G25").Value = "'20"
$ws.Range("G26").Value = "'20"
$ws.Range("D27").Value = "'0.1329"
$ws.Range("G27").Value = "'20"
$ws.Range("D28").Value = "'0.0002333"
$ws.Range("G28").Value = "'20"
$ws.Range("G29").Value = "'20"
$ws.Range("G30").Value = "'20"
$ws.Range("G31").Value = "'20"
$ws.Range("G32").Value = "'20"
$ws.Range("G33").Value = "'20"
$ws.Range("G34").Value = "'20"
$ws.Range("G35").Value = "'20"
$ws.Range("G36").Value = "'20"
$ws.Range("G37").Value = "'20"
$ws.Range("G38").Value = "'20"
$ws.Range("G39").Value = "'20"
$ws.Range("G40").Value = "'20"
$ws.Range("D41").Value = "'0.006264"
$ws.Range("E41").Value = "'40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'20"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("G42").Value = "'20"
$ws.Range("E43").Value = "'42CEJICEJI"
$ws.Range("G43").Value = "'20"
$ws.Range("D44").Value = "'0.008522"
$ws.Range("G44").Value = "'20"
$ws.Range("D45").Value = "'0.00005270"
$ws.Range("G45").Value = "'20"
$ws.Range("G46").Value = "'20"
$ws.Range("D47").Value = "'0.4403"
$ws.Range("G47").Value = "'20"
$ws.Range("D48").Value = "'0.002343"
$ws.Range("G48").Value = "'20"
$ws.Range("G49").Value = "'20"
$ws.Range("G50").Value = "'20"
$ws.Range("G51").Value = "'20"
